# Apply "Whatsapp web" registration-log update to the "Obras en general" sheet:
#   - add "Tasa de sellado" (col L) amounts for rows 81-84
#   - drop stray blank cells that used to live in row 88
#   - append five new work-log rows (89-93)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Obras en general")

# Plain assignment - used for values that Excel would never misinterpret as
# a number/date (so they keep behaving/looking exactly like ordinary text
# entered in this sheet, with no extra formatting).
function Set-Plain($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Text-forced assignment - prefixes the value with a leading apostrophe so
# Excel's automatic number/date detection can't turn numeric/date-looking
# text (e.g. "02/06/2025", "63305") into a real date/number, matching how
# the rest of the sheet stores everything as plain text. Also used (with an
# empty value) to create a real, present-but-empty cell in the saved XML
# instead of letting the cell be silently dropped - needed for the blank
# placeholder cells that sit alongside populated ones.
function Set-Text($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# ---------------------------------------------------------------------------
# Row 81 - add "Tasa de sellado" amount and the (blank) cells that sit
# alongside it once that column group is populated.
# ---------------------------------------------------------------------------
Set-Text "L81" "590000"
foreach ($col in "M","N","O","P","Q","T","U","V","W") {
    Set-Text "$col`81" ""
}

# ---------------------------------------------------------------------------
# Rows 82-84 - just add the "Tasa de sellado" amount.
# ---------------------------------------------------------------------------
Set-Text "L82" "31000"
Set-Text "L83" "55000"
Set-Text "L84" "26000"

# ---------------------------------------------------------------------------
# Row 88 - remove the stray blank cells that used to sit in D88, L88:Q88
# and T88:W88.
# ---------------------------------------------------------------------------
$ws.Range("D88").ClearContents()
$ws.Range("L88:Q88").ClearContents()
$ws.Range("T88:W88").ClearContents()

# ---------------------------------------------------------------------------
# New rows 89-93.
# ---------------------------------------------------------------------------
Set-Text  "A89" "02/06/2025"
Set-Plain "B89" "MMO"
Set-Plain "C89" "Físico"
Set-Plain "D89" "19 PLANOS"
Set-Plain "E89" "Obra nueva"
Set-Plain "F89" "CABRERA NATALIA PAMELA"
Set-Plain "G89" "CABALLERO ELIDA ROXANA"
Set-Plain "H89" "AV. UAN DOMINGO PERON (186) Nº 551"
Set-Plain "I89" "134/C/25"
Set-Text  "K89" "992"
Set-Plain "R89" "No pagado"
Set-Plain "S89" "No pagado"

Set-Text  "A90" "02/06/2025"
Set-Plain "B90" "MMO"
Set-Plain "C90" "Físico"
Set-Plain "D90" "19 PLANOS"
Set-Plain "E90" "Obra nueva"
Set-Plain "F90" "OTAZU JUAN JOSE "
Set-Plain "G90" "CABALLERO ELIDA ROXANA"
Set-Plain "H90" "AV. UAN DOMINGO PERON (186) Nº 551"
Set-Plain "I90" "134/C/25"
Set-Text  "K90" "992"
Set-Plain "L90" "None"
Set-Plain "M90" "None"
Set-Plain "N90" "None"
Set-Plain "O90" "None"
Set-Plain "P90" "None"
Set-Plain "Q90" "None"
Set-Plain "R90" "No pagado"
Set-Plain "S90" "No pagado"

Set-Text  "A91" "02/06/2025"
Set-Plain "B91" "MMO"
Set-Plain "C91" "Físico"
Set-Text  "D91" "16"
Set-Plain "E91" "Obra nueva"
Set-Plain "F91" "CABRERA NATALIA PAMELA"
Set-Plain "G91" "LOPEZ ENZO FACUNDO"
Set-Plain "H91" "CALLE PERU (104)- GARUPA"
Set-Plain "I91" "104/L/25"
Set-Text  "K91" "63305"
Set-Plain "R91" "No pagado"
Set-Plain "S91" "No pagado"

Set-Text  "A92" "02/06/2025"
Set-Plain "B92" "MMO"
Set-Plain "C92" "Físico"
Set-Text  "D92" "16"
Set-Plain "E92" "Obra nueva"
Set-Plain "F92" "OTAZU JUAN JOSE "
Set-Plain "G92" "LOPEZ ENZO FACUNDO"
Set-Plain "H92" "CALLE PERU (104)- GARUPA"
Set-Plain "I92" "104/L/25"
Set-Text  "K92" "63305"
Set-Plain "L92" "None"
Set-Plain "M92" "None"
Set-Plain "N92" "None"
Set-Plain "O92" "None"
Set-Plain "P92" "None"
Set-Plain "Q92" "None"
Set-Plain "R92" "No pagado"
Set-Plain "S92" "No pagado"

Set-Text  "A93" "02/06/2025"
Set-Plain "B93" "MMO"
Set-Plain "C93" "Físico"
Set-Text  "D93" "19"
Set-Plain "E93" "Registración"
Set-Plain "F93" "CABRERA NATALIA PAMELA"
Set-Plain "G93" "GOMEZ ROSSANA ELIZABETH"
Set-Plain "H93" "CALLE SALTO DEL MOCONA Nº 451- GARUPA"
Set-Plain "I93" "137/G/25"
Set-Text  "K93" "59140"
Set-Plain "R93" "No pagado"
Set-Plain "S93" "No pagado"

Write-Host "Applied registration log update (rows 81-93)."
